$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text of a run (found via its unique current text) with
# new text, without letting the engine auto-merge the edited run into an
# adjacent, identically-formatted neighbour run. The merge only happens when
# the edited span touches a run boundary, so we "fence" the target range by
# momentarily flipping Bold before the text write, then restore Bold (to its
# original value) on the freshly written range afterwards.
# ---------------------------------------------------------------------------
function Set-RunText($oldText, $newText) {
    $fr = $d.Content
    $found = $fr.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Set-RunText: could not find '$oldText'"
    }
    $start = $fr.Start
    $wasBold = $fr.Bold
    $fr.Bold = 1
    $fr.Text = $newText
    $r2 = $d.Range($start, $start + $newText.Length)
    $r2.Bold = $wasBold
    return $r2
}

# Append a brand-new run with the given text right after character offset
# $pos (copying the Aptos/black/24-half-pt body-text formatting used
# throughout this document), returning the new end offset.
function Add-Run($pos, $txt) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($txt)
    $newR = $d.Range($pos, $pos + $txt.Length)
    $newR.Font.NameAscii = "Aptos"
    $newR.Font.Name = "Aptos"
    $newR.Font.Size = 12
    $newR.Font.Color = 0
    return $pos + $txt.Length
}

# --- Title / byline / email ------------------------------------------------
Set-RunText "The Celestial Symphony: Unveiling the Harmony of the Cosmos" "Journeying Through the Realm of Chemistry: Unveiling the Secrets of Matter" | Out-Null
Set-RunText "Amelia Carter" "Emily Harper" | Out-Null
Set-RunText "ameliacarter@spaceobservatory" "emilyharper0622@ymail" | Out-Null
Set-RunText "edu" "net" | Out-Null

# --- Body paragraph 1 --------------------------------------------------
Set-RunText "In the vast expanse of the cosmos, celestial bodies dance to an intricate rhythm, creating a symphony of cosmic wonders" "In the vast expanse of human knowledge, chemistry stands as a captivating language that unlocks the hidden secrets of matter" | Out-Null

Set-RunText " The universe, like a celestial orchestra, weaves together the melodies of stars, planets, and galaxies, inviting us to decipher their harmonies and unravel the mysteries of existence" " It is a science that investigates the fundamental building blocks of the universe and the intricate interactions between them, guiding us towards an understanding of the world around us" | Out-Null

Set-RunText " From the blazing fires of the sun to the gentle glow of distant nebulas, each cosmic entity contributes to the grand composition, echoing the interconnectedness of all things" " From the smallest atoms to the sprawling galaxies, chemistry holds the key to unraveling the mysteries of existence" | Out-Null

Set-RunText "As we embark on this odyssey of cosmic exploration, we are granted a glimpse into the profound beauty and complexity of the universe" "Delving into the realm of chemistry, we embark on an extraordinary quest to explore the nature of substances and their transformations" | Out-Null

Set-RunText " We witness the birth and death of stars, the graceful ballet of planets around their suns, and the explosive spectacle of supernovae that herald the creation of new elements" " We learn about the elements, the basic units of matter, and the ways in which they combine to form compounds with diverse properties" | Out-Null

Set-RunText " With each observation, we deepen our understanding of the universe's composition, its evolution, and our place within this cosmic tapestry" " Through chemical reactions, we witness the dynamic interplay of atoms and molecules, as they rearrange and recombine, creating new substances with distinct characteristics" | Out-Null

Set-RunText "The study of celestial phenomena has ignited human curiosity for millennia, inspiring profound contemplations about our origins, our destiny, and the nature of reality itself" "Unraveling the enigmas of chemistry enables us to decipher the fundamental principles that govern the behavior of matter" | Out-Null

$r19 = Set-RunText " From ancient astronomers who charted the movements of celestial bodies to modern astrophysicists who probe the depths of space with cutting-edge telescopes, humanity's quest for knowledge about the cosmos has been an enduring endeavor, revealing the intricate workings of the universe and expanding our horizons of understanding" " We uncover the laws of thermodynamics, which dictate the flow of energy and the direction of chemical change"

# Insert four brand-new runs right after run 19 (before the trailing,
# unchanged "." run): ".", sentence, ".", sentence.
$pos = $r19.End
$pos = Add-Run $pos "."
$pos = Add-Run $pos " We unravel the mysteries of chemical bonding, the forces that hold atoms together and determine the properties of compounds"
$pos = Add-Run $pos "."
$pos = Add-Run $pos " Moreover, we delve into the intricacies of chemical reactions, exploring the factors that influence their rates and the mechanisms by which they occur"

# --- Summary paragraph ------------------------------------------------
Set-RunText "The celestial symphony that unfolds before us is a testament to the interconnectedness of the universe and the profound beauty that exists beyond our terrestrial sphere" "In this exploration of chemistry, we have embarked on a captivating journey through the realm of matter, unveiling the secrets of its structure, properties, and transformations" | Out-Null

Set-RunText " As we continue to unravel the mysteries of the cosmos, we gain a deeper appreciation for our place within the vastness of existence and the intricate harmonies that govern the symphony of the stars" " We have delved into the fundamental principles that govern chemical reactions and the intricate interactions between atoms and molecules" | Out-Null

# These two runs (the second of which carries a lastRenderedPageBreak) are
# replaced by a single new run.
Set-RunText " The study of celestial phenomena ignites our imagination, inspiring awe and wonder at the boundless marvels of the universe, reminding us that we are part of a grand cosmic dance that has been playing out for eons and will continue long after our own brief existence" " Through this exploration, we have gained a deeper understanding of the world around us and the remarkable complexity of the universe we inhabit" | Out-Null

# --- Trailing empty paragraph -------------------------------------------
$d.Paragraphs.Add() | Out-Null
